$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Successful run for 0.5s timestep: refresh the discharge-coefficient sweep
# results (columns B:I, rows 2-11) with the values produced by the new run.
$newData = @{
    2  = @{ B = 33.285905779941892;  C = 16.275184844542355;  D = 0.48895123816488695; E = 30.43706377061801;  F = 15.354937458636416; G = 0.50448156150525558; H = 281.5; I = 247.5 }
    3  = @{ B = 33.316077996955912;  C = 16.303421194807829;  D = 0.48935595589305175; E = 30.501523533415934; F = 15.415012223537902; G = 0.50538499188901143; H = 281.5; I = 247.5 }
    4  = @{ B = 33.341114017191778;  C = 16.328640322741769;  D = 0.48974489317670022; E = 30.494467243112851; F = 15.397149174403538; G = 0.50491615582761062; H = 281.5; I = 247.5 }
    5  = @{ B = 33.361589039801849;  C = 16.348184105031642;  D = 0.49003013871813889; E = 30.515752826899742; F = 15.399669414184094; G = 0.5046465509646294;  H = 281.5; I = 247.5 }
    6  = @{ B = 33.376007316432542;  C = 16.307877027489013;  D = 0.48861078177736067; E = 30.545213126406018; F = 15.415925728397664; G = 0.50469203356354242; H = 281.5; I = 247   }
    7  = @{ B = 33.389708971335352;  C = 16.311471467215917;  D = 0.4885179287192683;  E = 30.587524235588717; F = 15.482524661983135; G = 0.50617122663262659; H = 281.5; I = 247   }
    8  = @{ B = 33.399523333405867;  C = 16.322743477919445;  D = 0.48871186917789339; E = 30.602007772780688; F = 15.494615353011225; G = 0.50632675699118967; H = 281.5; I = 247   }
    9  = @{ B = 33.409746698148012;  C = 16.329580341143316;  D = 0.48876696039268402; E = 30.592300539390152; F = 15.488203735858297; G = 0.50627783667056803; H = 281.5; I = 246.5 }
    10 = @{ B = 33.419073408951213;  C = 16.342886520766751;  D = 0.48902871485327748; E = 30.615062913747778; F = 15.476099470649883; G = 0.50550604825640577; H = 281.5; I = 246.5 }
    11 = @{ B = 33.428031668027195;  C = 16.348906240676111;  D = 0.48907774178978353; E = 30.61591614602488;  F = 15.457927257041952; G = 0.50489840589170099; H = 281.5; I = 246.5 }
}

foreach ($row in $newData.Keys) {
    $rowValues = $newData[$row]
    foreach ($col in $rowValues.Keys) {
        $ws.Range("$col$row").Value = $rowValues[$col]
    }
}

# The wider numbers no longer fit the old column widths, so the columns were
# widened (same action AutoFit would perform) to comfortably fit the new data.
$newColumnWidths = @{
    1 = 24.166666666666668
    2 = 32.5
    3 = 31.0
    4 = 37.833333333333336
    5 = 31.166666666666668
    6 = 29.666666666666668
    7 = 36.666666666666664
    8 = 31.5
    9 = 30.166666666666668
}

foreach ($col in $newColumnWidths.Keys) {
    $ws.Columns.Item($col).ColumnWidth = $newColumnWidths[$col]
}
